$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Word maintains a single "_GoBack" bookmark marking the last edit spot.
#    In the starting document it sits near the "Final Project" heading; the
#    edit we are replaying relocates it to the newly authored list item
#    below, so drop the old one before we insert the new one (a document can
#    only have one bookmark with a given name).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# 2. Locate the existing "(OpenIntro) Chapter 8..." reading-list bullet and
#    add a sibling bullet right after it for the new OpenIntro chapter link.
# ---------------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute("(OpenIntro) Chapter 8: Introduction to linear regression", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $findRange.Find.Found) {
    throw "Could not locate the Chapter 8 reading list paragraph"
}

$chapter8Para = $findRange.Paragraphs(1)
$chapter8Index = $chapter8Para.Index

# InsertParagraphAfter clones the paragraph/run formatting (numbering, tabs,
# spacing, indent, fonts) from the Chapter 8 bullet, so the new bullet lines
# up with the rest of the list.
$chapter8Para.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs($chapter8Index + 1)
$newParaRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

# Replace the (empty, cloned-formatting) paragraph's contents with the exact
# target markup: list paragraph properties, the relocated "_GoBack" bookmark,
# and a run (without the inherited black font color) holding the new text.
$newParaXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="44"/></w:numPr><w:tabs><w:tab w:val="clear" w:pos="720"/><w:tab w:val="num" w:pos="1152"/></w:tabs><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:ind w:left="1152"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>(OpenIntro) Chapter 9: Multiple and Logistic Regression</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newParaRange.InsertXML($newParaXml) | Out-Null

Write-Output "Inserted OpenIntro Chapter 9 reading list item and relocated _GoBack bookmark"
